# Updated: po 08. 02. 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrections to existing AgTests (H) / AgPosit (I) figures ---
$ws.Range("H313").Value = 61379

$ws.Range("H321").Value = 90731
$ws.Range("I321").Value = 2794

$ws.Range("H322").Value = 107036
$ws.Range("I322").Value = 2304

$ws.Range("H323").Value = 148636
$ws.Range("I323").Value = 2284

$ws.Range("H324").Value = 231232
$ws.Range("I324").Value = 2653

$ws.Range("H325").Value = 723880
$ws.Range("I325").Value = 5986

$ws.Range("H326").Value = 425220
$ws.Range("I326").Value = 3737

$ws.Range("H327").Value = 235963
$ws.Range("I327").Value = 2874

$ws.Range("H328").Value = 178431
$ws.Range("I328").Value = 2618

$ws.Range("H329").Value = 82218

$ws.Range("H330").Value = 70830
$ws.Range("I330").Value = 1988

$ws.Range("H331").Value = 148614
$ws.Range("I331").Value = 2567

$ws.Range("H332").Value = 415219
$ws.Range("I332").Value = 4051

$ws.Range("H333").Value = 255932
$ws.Range("I333").Value = 2740

$ws.Range("H334").Value = 202337
$ws.Range("I334").Value = 3406

$ws.Range("H335").Value = 124096
$ws.Range("I335").Value = 2864

$ws.Range("H336").Value = 100009
$ws.Range("I336").Value = 3155

$ws.Range("H337").Value = 102098
$ws.Range("I337").Value = 2875

# --- Append three new days of data (rows 338-340) ---
$newRows = @(
    @{ Row = 338; A = 44232; B = 261774; C = 238221; D = 18418; E = 13677; F = 2241; G = 5135; H = 211343; I = 3726 },
    @{ Row = 339; A = 44233; B = 263326; C = 241686; D = 16441; E = 8282;  F = 1552; G = 5199; H = 511033; I = 4671 },
    @{ Row = 340; A = 44234; B = 264083; C = 245702; D = 13110; E = 3784;  F = 757;  G = 5271; H = 284609; I = 2471 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}
